# Atualização de bases das ligas, do dia: 11-06-2024 às 21:19
#
# The source data had several fixtures whose match-result/odds rows were
# attributed to the wrong fixture. This script swaps the data (id/odds/
# score columns) between the affected row pairs while leaving the row's
# own "id" (column A) and "Date" (column D) in place, matching the
# canonical fix. For the very first pair (rows 3/4) the HomeTeam names
# were also corrected directly (their underlying labels were swapped),
# so the HomeTeam column is restored afterwards instead of being swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRanges {
    param($ws, $row1, $row2, $firstCol, $lastCol)

    $rangeA = $ws.Range($firstCol + $row1 + ":" + $lastCol + $row1)
    $rangeB = $ws.Range($firstCol + $row2 + ":" + $lastCol + $row2)

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# Rows 3 and 4: swap everything from column B through AD (ids, odds, score, etc.)
Swap-RowRanges $ws 3 4 "B" "AD"
# HomeTeam (column E) must keep showing the same two team names as before
# (the fix here was purely to the match-data columns, not the team labels).
$ws.Range("E3").Value = "Magdeburg II"
$ws.Range("E4").Value = "SV Dessau 05"

# Rows 86 and 87: full swap of B through AD.
Swap-RowRanges $ws 86 87 "B" "AD"

# Rows 120 and 121: full swap of B through AD.
Swap-RowRanges $ws 120 121 "B" "AD"

# Rows 151 and 152: full swap of B through AD.
Swap-RowRanges $ws 151 152 "B" "AD"

# Rows 180 and 181: full swap of B through AD.
Swap-RowRanges $ws 180 181 "B" "AD"
